$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 311, shifting existing rows 311-318 down to 312-319
$ws.Rows.Item(311).Insert()

# Populate the newly inserted row 311 with the new data record
$ws.Range("A311").Value = 11
$ws.Range("B311").Value = "Vega Monumental Concepción"
$ws.Range("C311").Value = "Bíobío"
$ws.Range("D311").Value = 45239
$ws.Range("D311").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E311").Value = 8
$ws.Range("F311").Value = "Fruta"
$ws.Range("G311").Value = 100108
$ws.Range("H311").Value = "Tropicales y subtropicales"
$ws.Range("I311").Value = 100108005
$ws.Range("J311").Value = "Piña"
$ws.Range("K311").Value = "Caramelo"
$ws.Range("L311").Value = "Segunda"
$ws.Range("M311").Value = 250
$ws.Range("N311").Value = 22000
$ws.Range("O311").Value = 23000
$ws.Range("P311").Value = 22400
$ws.Range("Q311").Value = "$/caja 14 unidades"
$ws.Range("R311").Value = "Ecuador"
$ws.Range("S311").Value = 1600
$ws.Range("T311").Value = 14
